$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: names (write first so new shared strings land in this order) ---
$ws.Range("A2").Value = "Eric Bradshaw"
$ws.Range("A3").Value = "Eric Wong"
$ws.Range("A4").Value = "Eliass Ghauss"

# --- Column B: emails ---
$ws.Range("B2").Value = "ericbradshaw@gmail.com"
$ws.Range("B3").Value = "ericwong@gmail.com"
$ws.Range("B4").Value = "eliassghauss@gmail.com"

# --- Column C: class (new row 4) ---
$ws.Range("C4").Value = "B"

# --- Column D: password hash (new row 4) ---
$ws.Range("D4").Value = "`$2y`$12`$JNv6zD/l1zyyeJ6S8aMSPeUkg9.ZPB7qyOAGL8GavfXPuTYzjRsTa"

# --- Column E: role ---
$ws.Range("E3").Value = "teacher"
$ws.Range("E4").Value = "employee"

# --- Column F: status (new row 4) ---
$ws.Range("F4").Value = "searching"

# --- Column G: description (new row 4) ---
$ws.Range("G4").Value = "testing"

# --- Column H: school_id (new row 4) ---
$ws.Range("H4").Value = 1

# --- Columns I/J: profile_image / company_name (new row 4) ---
$ws.Range("I4").Value = "null"
$ws.Range("J4").Value = "null"

# --- Columns K/L/M/N: company / working / interviewing / searching flags (new row 4) ---
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0

# --- Column O: skills -- cleared out for rows 2 and 3 (row 2 keeps an empty styled cell,
#     row 3's cell is removed entirely, row 4 never gets one) ---
$ws.Range("O2").ClearContents()
$ws.Range("O3").ClearContents()

# --- Column P: positiontitle ---
$ws.Range("P3").Value = "Teacher"
$ws.Range("P4").Value = "Hiring Manager"

# --- Hyperlinks: rebuild mailto links for the (now changed) B column addresses ---
$origB2Style = $ws.Range("B2").Style
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:ericbradshaw@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:ericwong@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:eliassghauss@gmail.com")

# Restore the original (Courier New) hyperlink-like style on B2/B3, and match it on the new B4
$ws.Range("B2").Style = $origB2Style
$ws.Range("B3").Style = $origB2Style
$ws.Range("B4").Style = $origB2Style

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 12.1666666666667
$ws.Columns.Item(2).ColumnWidth = 22.1666666666667

# --- Selection ---
$ws.Range("O4").Select()
